$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OKE")

# Row 4 - Inventory
$ws.Range("B4").Value = 513000000.0
$ws.Range("C4").Value = 371000000.0
$ws.Range("D4").Value = 460000000.0
$ws.Range("E4").Value = 515000000.0
$ws.Range("F4").Value = 386000000.0

# Row 15 - Accounts Payable
$ws.Range("B15").Value = 930000000.0
$ws.Range("C15").Value = 719000000.0
$ws.Range("D15").Value = 624000000.0
$ws.Range("E15").Value = 755000000.0
$ws.Range("F15").Value = 742000000.0

# Row 24 - Long Term Tax Liability (Deferred)
$ws.Range("B24").Value = 815000000.0
$ws.Range("C24").Value = 670000000.0
$ws.Range("D24").Value = 582000000.0
$ws.Range("E24").Value = 476000000.0
$ws.Range("F24").Value = 443000000.0

# Row 38 - Net Debt
$ws.Range("G38").Value = 12466452000.0

# Row 39 - Total Debt
$ws.Range("G39").Value = 12487410000.0
